$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new greeting value to Sheet1!B4 (shared string), matching the
# centered alignment style already used by the other data cells.
$ws1.Range("B4").Value = "नमस्ते"
$ws1.Range("B4").HorizontalAlignment = $ws1.Range("A1").HorizontalAlignment

# Define a new named range "Χαιρετισμός" pointing at Sheet1!$B$4.
# (Creating a name whose first character is non-ASCII directly via
# Names.Add / Range.Name fails in this runtime, so create it with a
# plain placeholder name first and then rename it.)
$ws1.Range("B4").Name = "GreetingNamePlaceholder"
$newName = $wb.Names.Item("GreetingNamePlaceholder")
$newName.Name = "Χαιρετισμός"

# Select B4 on Sheet1 to match the saved selection/active-cell state.
$ws1.Activate()
$ws1.Range("B4").Select()

$wb.Save()
